$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.27782154083252
$ws.Range("B1").Value = 2.46479606628418
$ws.Range("C1").Value = 3.409960985183716
$ws.Range("D1").Value = 3.214523553848267
$ws.Range("E1").Value = 1.066484332084656
